# Weekly update: two new price records for "Vega Modelo de Temuco - Brócoli"
# are inserted ahead of the existing row 283, pushing the rest of the table
# (old rows 283-300) down by two rows (to 285-302).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 283/284; Excel shifts old rows 283-300 down to 285-302
# and carries the formatting (incl. the date style on column D) down from row 282.
$ws.Range("A283:A284").EntireRow.Insert()

# --- New row 283 ---
$ws.Cells.Item(283, 1).Value = 10
$ws.Cells.Item(283, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(283, 3).Value = "La Araucanía"
$ws.Cells.Item(283, 4).Value = 44516
$ws.Cells.Item(283, 5).Value = 9
$ws.Cells.Item(283, 6).Value = 100112023
$ws.Cells.Item(283, 7).Value = "Brócoli"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Primera"
$ws.Cells.Item(283, 10).Value = 1250
$ws.Cells.Item(283, 11).Value = 800
$ws.Cells.Item(283, 12).Value = 800
$ws.Cells.Item(283, 13).Value = 800
$ws.Cells.Item(283, 14).Value = "$/unidad"
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 800
$ws.Cells.Item(283, 17).Value = 1
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# --- New row 284 ---
$ws.Cells.Item(284, 1).Value = 10
$ws.Cells.Item(284, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(284, 3).Value = "La Araucanía"
$ws.Cells.Item(284, 4).Value = 44516
$ws.Cells.Item(284, 5).Value = 9
$ws.Cells.Item(284, 6).Value = 100112023
$ws.Cells.Item(284, 7).Value = "Brócoli"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 2500
$ws.Cells.Item(284, 11).Value = 800
$ws.Cells.Item(284, 12).Value = 800
$ws.Cells.Item(284, 13).Value = 800
$ws.Cells.Item(284, 14).Value = "$/unidad"
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 800
$ws.Cells.Item(284, 17).Value = 1
$ws.Cells.Item(284, 18).Value = "Hortaliza"
